# Preliminary updates for neighbourhood-code refactor of the interface.
#
# The template's first three columns used to be:
#   A = districtcode, B = districtname, C = neighbourhoodcode
# They are refactored to:
#   A = neighbourhoodcode, B = neighbourhoodname, C = districtcode
# ("districtname" is retired and a brand-new "neighbourhoodname" string
# takes its place in the shared-strings table.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header relabel (row 1) ---------------------------------------------
$ws.Range("A1").Value = "neighbourhoodcode"
$ws.Range("B1").Value = "neighbourhoodname"
$ws.Range("C1").Value = "districtcode"

# --- Column widths (A:D) re-tightened to fit the new header text --------
# (closest values reachable through the ColumnWidth grid of this engine)
$ws.Columns.Item(1).ColumnWidth = 10.6666666667
$ws.Columns.Item(2).ColumnWidth = 19.3333333333
$ws.Columns.Item(3).ColumnWidth = 9.5
$ws.Columns.Item(4).ColumnWidth = 14.8333333333

# --- Active selection moved from G7 to F11 -------------------------------
[void]$ws.Range("F11").Select()
